# Set the "Industries" column (H) values to 0 for rows 27 through 73.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H27:H73").Value = 0
